# Updates cryptos list (Price / Volume(1h) columns) to match the
# latest scrape, including a row re-ordering where "Dai" and "PEPE"
# swap places (rows 38 and 39).
#
# Several "Price" values look like plain decimal numbers (e.g. "0.999",
# "18.03"). Excel's Range.Value setter auto-detects such strings and
# stores them as numeric cells, which would lose the original text
# formatting (trailing zeros, cell type). To keep them as text - as in
# the source workbook - we briefly force a text number format before
# assigning the value, then restore the cell's style so no stray
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.691.93"
$ws.Range("E2").Value = "  +9.29%  "
$ws.Range("D3").Value = "3.476.64"
$ws.Range("E3").Value = "  +13.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +13.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "545.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.63%  "
$ws.Range("D7").Value = "3.472.18"
$ws.Range("E7").Value = "  +13.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.627"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +20.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000266"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.63%  "
$ws.Range("D15").Value = "4.020.51"
$ws.Range("E15").Value = "  +12.47%  "
$ws.Range("D16").Value = "3.460.37"
$ws.Range("E16").Value = "  +12.26%  "
$ws.Range("E17").Value = "  +8.89%  "
$ws.Range("D18").Value = "66.612.29"
$ws.Range("E18").Value = "  +9.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.986"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "416.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +18.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "651.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.110"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.37%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0811"
$ws.Range("E38").Value = "  +24.99%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.388"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.137"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +20.12%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "2.987.21"
$ws.Range("E44").Value = "  +8.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0413"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +21.34%  "
$ws.Range("E51").Value = "  +9.80%  "
